$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '54.288.52'
$ws.Range("E2").Value = '  +0.61%  '

$ws.Range("D3").Value = '2.266.89'
$ws.Range("E3").Value = '  +1.20%  '

$ws.Range("E4").Value = '  -0.54%  '

$c = $ws.Range("D5")
$c.Value = "'497.09"
$c.Style = "Normal"
$ws.Range("E5").Value = '  +0.72%  '

$c = $ws.Range("D6")
$c.Value = "'128.97"
$c.Style = "Normal"
$ws.Range("E6").Value = '  +1.24%  '

$ws.Range("E7").Value = '  -0.18%  '

$ws.Range("E8").Value = '  -0.25%  '

$c = $ws.Range("D9")
$c.Value = "'0.0953"
$c.Style = "Normal"
$ws.Range("E9").Value = '  +0.57%  '

$ws.Range("E10").Value = '  +0.96%  '

$c = $ws.Range("D11")
$c.Value = "'0.335"
$c.Style = "Normal"
$ws.Range("E11").Value = '  +2.95%  '

$c = $ws.Range("D12")
$c.Value = "'4.89"
$c.Style = "Normal"
$ws.Range("E12").Value = '  +5.44%  '

$ws.Range("E13").Value = '  +5.83%  '

$ws.Range("D14").Value = '2.667.29'
$ws.Range("E14").Value = '  -0.04%  '

$ws.Range("D15").Value = '54.251.11'
$ws.Range("E15").Value = '  +0.25%  '

$ws.Range("E16").Value = '  +0.73%  '

$ws.Range("D17").Value = '2.271.18'
$ws.Range("E17").Value = '  -0.26%  '

$c = $ws.Range("D18")
$c.Value = "'10.24"
$c.Style = "Normal"
$ws.Range("E18").Value = '  +2.49%  '

$c = $ws.Range("D19")
$c.Value = "'4.14"
$c.Style = "Normal"
$ws.Range("E19").Value = '  +1.25%  '

$c = $ws.Range("D20")
$c.Value = "'303.85"
$c.Style = "Normal"
$ws.Range("E20").Value = '  +1.41%  '

$ws.Range("E21").Value = '  -1.26%  '

$ws.Range("E22").Value = '  +0.21%  '

$c = $ws.Range("D23")
$c.Value = "'60.62"
$c.Style = "Normal"
$ws.Range("E23").Value = '  -2.47%  '

$c = $ws.Range("D24")
$c.Value = "'1.00"
$c.Style = "Normal"
$ws.Range("E24").Value = '  +0.02%  '

$ws.Range("E25").Value = '  +0.71%  '

$c = $ws.Range("D26")
$c.Value = "'7.31"
$c.Style = "Normal"
$ws.Range("E26").Value = '  +4.07%  '

$c = $ws.Range("D27")
$c.Value = "'173.24"
$c.Style = "Normal"
$ws.Range("E27").Value = '  +3.86%  '

$ws.Range("B28").Value = 'PEPE'
$ws.Range("C28").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D28").Value = '0.0₃0693'
$ws.Range("E28").Value = '  +1.58%  '

$ws.Range("B29").Value = 'PancakeSwap'
$ws.Range("C29").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$c = $ws.Range("D29")
$c.Value = "'1.60"
$c.Style = "Normal"
$ws.Range("E29").Value = '  +0.57%  '

$c = $ws.Range("D30")
$c.Value = "'5.97"
$c.Style = "Normal"
$ws.Range("E30").Value = '  +2.23%  '

$ws.Range("E31").Value = '  +1.53%  '

$ws.Range("E32").Value = '  +0.00%  '

$c = $ws.Range("D33")
$c.Value = "'17.82"
$c.Style = "Normal"
$ws.Range("E33").Value = '  +0.95%  '

$ws.Range("E34").Value = '  +0.26%  '

$c = $ws.Range("D35")
$c.Value = "'0.948"
$c.Style = "Normal"
$ws.Range("E35").Value = '  +5.70%  '

$c = $ws.Range("D36")
$c.Value = "'1.21"
$c.Style = "Normal"
$ws.Range("E36").Value = '  +1.95%  '

$c = $ws.Range("D37")
$c.Value = "'3.71"
$c.Style = "Normal"
$ws.Range("E37").Value = '  +1.04%  '

$ws.Range("E38").Value = '  +0.84%  '

$ws.Range("E39").Value = '  +0.32%  '

$ws.Range("E40").Value = '  +0.93%  '

$c = $ws.Range("D41")
$c.Value = "'4.81"
$c.Style = "Normal"
$ws.Range("E41").Value = '  +1.74%  '

$c = $ws.Range("D42")
$c.Value = "'124.76"
$c.Style = "Normal"
$ws.Range("E42").Value = '  -0.74%  '

$c = $ws.Range("D43")
$c.Value = "'0.0492"
$c.Style = "Normal"
$ws.Range("E43").Value = '  +2.12%  '

$c = $ws.Range("D44")
$c.Value = "'0.0895"
$c.Style = "Normal"
$ws.Range("E44").Value = '  +1.38%  '

$c = $ws.Range("D45")
$c.Value = "'0.546"
$c.Style = "Normal"
$ws.Range("E45").Value = '  +0.67%  '

$c = $ws.Range("D46")
$c.Value = "'242.55"
$c.Style = "Normal"
$ws.Range("E46").Value = '  +2.73%  '

$ws.Range("E47").Value = '  +1.10%  '

$c = $ws.Range("D48")
$c.Value = "'0.0205"
$c.Style = "Normal"
$ws.Range("E48").Value = '  +1.61%  '

$ws.Range("E49").Value = '  +0.76%  '

$c = $ws.Range("D50")
$c.Value = "'16.17"
$c.Style = "Normal"
$ws.Range("E50").Value = '  +0.42%  '

$c = $ws.Range("D51")
$c.Value = "'0.936"
$c.Style = "Normal"
$ws.Range("E51").Value = '  -0.12%  '
